# Natmi following Dr Hou advice
# Updates the Gnai2-Tshr LR-pair sheet: existing rows 2-13 get refreshed
# NATMI metrics, and 4 new rows (14-17) are appended for the sCs sending
# cluster against all four target clusters (ECs, FAPs, M2, sCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-17 (A..T) reflecting the updated NATMI Gnai2-Tshr LR-pair analysis
$data = @(
    @("ECs", "Gnai2", "Tshr", "ECs", 3, 1, 162.399297, 487.197891, 0.3910371682630009, 0.3910371682630009, 3, 1, 0.7209850000000001, 2.162955, 0.1533325535399077, 0.1533325535399077, 117.087457147545, 1053.787114327905, 0.0599587275387805, 0.0599587275387805),
    @("ECs", "Gnai2", "Tshr", "FAPs", 3, 1, 162.399297, 487.197891, 0.3910371682630009, 0.3910371682630009, 3, 1, 3.387303666666666, 10.161911, 0.7203810354238886, 0.7203810354238885, 550.095734192189, 4950.861607729701, 0.281695760162526, 0.281695760162526),
    @("ECs", "Gnai2", "Tshr", "M2", 3, 1, 162.399297, 487.197891, 0.3910371682630009, 0.3910371682630009, 2, 0.6666666666666666, 0.1167496666666667, 0.350249, 0.02482926068494219, 0.02482926068494219, 18.960063791651, 170.640574124859, 0.009709163788303655, 0.009709163788303653),
    @("ECs", "Gnai2", "Tshr", "sCs", 3, 1, 162.399297, 487.197891, 0.3910371682630009, 0.3910371682630009, 3, 1, 0.4770616666666667, 1.431185, 0.1014571503512615, 0.1014571503512615, 77.47447929231501, 697.2703136308351, 0.03967351677339084, 0.03967351677339083),
    @("FAPs", "Gnai2", "Tshr", "ECs", 3, 1, 65.41736466666667, 196.252094, 0.1575168212364948, 0.1575168212364948, 3, 1, 0.7209850000000001, 2.162955, 0.1533325535399077, 0.1533325535399077, 47.16493866419668, 424.4844479777701, 0.02415245642568092, 0.02415245642568091),
    @("FAPs", "Gnai2", "Tshr", "FAPs", 3, 1, 65.41736466666667, 196.252094, 0.1575168212364948, 0.1575168212364948, 3, 1, 3.387303666666666, 10.161911, 0.7203810354238886, 0.7203810354238885, 221.5884791990704, 1994.296312791634, 0.1134721307790257, 0.1134721307790257),
    @("FAPs", "Gnai2", "Tshr", "M2", 3, 1, 65.41736466666667, 196.252094, 0.1575168212364948, 0.1575168212364948, 2, 0.6666666666666666, 0.1167496666666667, 0.350249, 0.02482926068494219, 0.02482926068494219, 7.637455519045113, 68.737099671406, 0.003911026216744369, 0.003911026216744368),
    @("FAPs", "Gnai2", "Tshr", "sCs", 3, 1, 65.41736466666667, 196.252094, 0.1575168212364948, 0.1575168212364948, 3, 1, 0.4770616666666667, 1.431185, 0.1014571503512615, 0.1014571503512615, 31.20811701682112, 280.87305315139, 0.01598120781504384, 0.01598120781504384),
    @("M2", "Gnai2", "Tshr", "ECs", 3, 1, 126.3069433333333, 378.92083, 0.3041313008456065, 0.3041313008456065, 3, 1, 0.7209850000000001, 2.162955, 0.1533325535399077, 0.1533325535399077, 91.06541153918334, 819.5887038526502, 0.04663322897007074, 0.04663322897007074),
    @("M2", "Gnai2", "Tshr", "FAPs", 3, 1, 126.3069433333333, 378.92083, 0.3041313008456065, 0.3041313008456065, 3, 1, 3.387303666666666, 10.161911, 0.7203810354238886, 0.7203810354238885, 427.8399722784589, 3850.55975050613, 0.2190904214079722, 0.2190904214079722),
    @("M2", "Gnai2", "Tshr", "M2", 3, 1, 126.3069433333333, 378.92083, 0.3041313008456065, 0.3041313008456065, 2, 0.6666666666666666, 0.1167496666666667, 0.350249, 0.02482926068494219, 0.02482926068494219, 14.74629353185222, 132.71664178667, 0.007551355351146143, 0.007551355351146143),
    @("M2", "Gnai2", "Tshr", "sCs", 3, 1, 126.3069433333333, 378.92083, 0.3041313008456065, 0.3041313008456065, 3, 1, 0.4770616666666667, 1.431185, 0.1014571503512615, 0.1014571503512615, 60.25620089817223, 542.3058080835501, 0.03085629511641744, 0.03085629511641744),
    @("sCs", "Gnai2", "Tshr", "ECs", 3, 1, 61.180387, 183.541161, 0.1473147096548978, 0.1473147096548978, 3, 1, 0.7209850000000001, 2.162955, 0.1533325535399077, 0.1533325535399077, 44.110141321195, 396.991271890755, 0.02258814060537558, 0.02258814060537558),
    @("sCs", "Gnai2", "Tshr", "FAPs", 3, 1, 61.180387, 183.541161, 0.1473147096548978, 0.1473147096548978, 3, 1, 3.387303666666666, 10.161911, 0.7203810354238886, 0.7203810354238885, 207.2365492131856, 1865.128942918671, 0.1061227230743648, 0.1061227230743648),
    @("sCs", "Gnai2", "Tshr", "M2", 3, 1, 61.180387, 183.541161, 0.1473147096548978, 0.1473147096548978, 2, 0.6666666666666666, 0.1167496666666667, 0.350249, 0.02482926068494219, 0.02482926068494219, 7.142789788787667, 64.285108099089, 0.003657715328748028, 0.003657715328748028),
    @("sCs", "Gnai2", "Tshr", "sCs", 3, 1, 61.180387, 183.541161, 0.1473147096548978, 0.1473147096548978, 3, 1, 0.4770616666666667, 1.431185, 0.1014571503512615, 0.1014571503512615, 29.18681738953167, 262.681356505785, 0.0149461306464094, 0.0149461306464094),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
